# Update cryptos list: apply per-cell value changes from the commit diff.
# Note: a few "Price" cells are text like "7.50"/"1.00"/"6.40"/"0.130"/"31.30"
# whose trailing zero(s) would be silently stripped if Excel auto-coerced the
# assigned string into a number. Prefixing those specific values with a
# leading apostrophe forces Excel to keep them as literal text (same as the
# other already-text Price cells), without touching any cell's NumberFormat.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '43.222.80'
$ws.Range('E2').Value = '  +5.02%  '

$ws.Range('D3').Value = '2.303.77'
$ws.Range('E3').Value = '  +5.74%  '

$ws.Range('E4').Value = '  +0.11%  '

$ws.Range('D5').Value = '251.67'
$ws.Range('E5').Value = '  +0.63%  '

$ws.Range('D6').Value = '0.643'
$ws.Range('E6').Value = '  +5.04%  '

$ws.Range('D7').Value = '73.65'
$ws.Range('E7').Value = '  +10.93%  '

$ws.Range('E8').Value = '  -0.02%  '

$ws.Range('D9').Value = '0.657'
$ws.Range('E9').Value = '  +11.10%  '

$ws.Range('D10').Value = '39.59'
$ws.Range('E10').Value = '  +6.58%  '

$ws.Range('D11').Value = '0.0984'
$ws.Range('E11').Value = '  +5.29%  '

$ws.Range('D12').Value = '59.68'
$ws.Range('E12').Value = '  +0.84%  '

$ws.Range('D13').Value = "'7.50"
$ws.Range('E13').Value = '  +9.19%  '

$ws.Range('D14').Value = '0.105'
$ws.Range('E14').Value = '  +1.55%  '

$ws.Range('D15').Value = '2.648.61'
$ws.Range('E15').Value = '  +5.62%  '

$ws.Range('D16').Value = '15.28'
$ws.Range('E16').Value = '  +7.22%  '

$ws.Range('D17').Value = '0.897'
$ws.Range('E17').Value = '  +5.93%  '

$ws.Range('D18').Value = '2.309.16'
$ws.Range('E18').Value = '  +6.62%  '

$ws.Range('D19').Value = '43.235.96'
$ws.Range('E19').Value = '  +5.09%  '

$ws.Range('E20').Value = '  +6.92%  '

$ws.Range('D21').Value = '6.41'
$ws.Range('E21').Value = '  +6.11%  '

$ws.Range('E22').Value = '  +3.11%  '

$ws.Range('B23').Value = 'ImmutableX'
$ws.Range('C23').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D23').Value = '2.29'
$ws.Range('E23').Value = '  +13.87%  '

$ws.Range('B24').Value = 'BitcoinCash'
$ws.Range('C24').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D24').Value = '236.71'
$ws.Range('E24').Value = '  +2.99%  '

$ws.Range('E25').Value = '  +2.79%  '

$ws.Range('D26').Value = '11.85'
$ws.Range('E26').Value = '  +5.90%  '

$ws.Range('D27').Value = "'1.00"
$ws.Range('E27').Value = '  +0.05%  '

$ws.Range('E28').Value = '  +2.30%  '

$ws.Range('E29').Value = '  -0.09%  '

$ws.Range('D30').Value = '2.18'
$ws.Range('E30').Value = '  +7.12%  '

$ws.Range('D31').Value = '168.17'
$ws.Range('E31').Value = '  +0.65%  '

$ws.Range('D32').Value = '21.38'
$ws.Range('E32').Value = '  +6.16%  '

$ws.Range('D33').Value = "'6.40"
$ws.Range('E33').Value = '  +12.06%  '

$ws.Range('D34').Value = "'0.130"
$ws.Range('E34').Value = '  +8.20%  '

$ws.Range('D35').Value = '0.0816'
$ws.Range('E35').Value = '  +8.76%  '

$ws.Range('D36').Value = "'31.30"
$ws.Range('E36').Value = '  +25.17%  '

$ws.Range('D37').Value = '0.126'
$ws.Range('E37').Value = '  +3.87%  '

$ws.Range('D38').Value = '4.74'
$ws.Range('E38').Value = '  +16.77%  '

$ws.Range('D39').Value = '4.83'
$ws.Range('E39').Value = '  +6.96%  '

$ws.Range('D40').Value = '0.0313'
$ws.Range('E40').Value = '  +2.62%  '

$ws.Range('D41').Value = '13.75'
$ws.Range('E41').Value = '  +20.56%  '

$ws.Range('D42').Value = '2.37'
$ws.Range('E42').Value = '  +7.17%  '

$ws.Range('D43').Value = '6.13'
$ws.Range('E43').Value = '  +10.90%  '

$ws.Range('D44').Value = '0.214'
$ws.Range('E44').Value = '  +12.59%  '

$ws.Range('D45').Value = '9.29'
$ws.Range('E45').Value = '  +9.08%  '

$ws.Range('D46').Value = '63.11'
$ws.Range('E46').Value = '  +4.23%  '

$ws.Range('D47').Value = '4.97'
$ws.Range('E47').Value = '  -8.93%  '

$ws.Range('E48').Value = '  +6.11%  '

$ws.Range('E49').Value = '  +5.63%  '

$ws.Range('E50').Value = '  +0.18%  '

$ws.Range('E51').Value = '  +5.43%  '
